# reporter: Use the term "issue" instead of "error"
#
# The scan report workbook has two sheets ("Summary" and the per-project
# sheet) that both carry the same header row (row 10) with the column
# headers "Analyzer Errors" / "Scan Errors" in columns E/F. Rename those
# headers to "Analyzer Issues" / "Scan Issues" on every sheet, and restore
# the (now slightly different, because the new text has a different
# best-fit pixel width) auto-fit width of column E.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    if ($ws.Range("E10").Value2 -eq "Analyzer Errors") {
        $ws.Range("E10").Value = "Analyzer Issues"
    }
    if ($ws.Range("F10").Value2 -eq "Scan Errors") {
        $ws.Range("F10").Value = "Scan Issues"
    }

    # Column E was previously best-fit to "Analyzer Errors" (14.9765625
    # characters wide); re-fit it to the new "Analyzer Issues" text, which
    # is fractionally wider (15.00390625 characters, i.e. raw OOXML
    # <col> width). The host's ColumnWidth setter re-derives the stored
    # raw width as round(6 * ColumnWidth + 5) / 6, so feed it the value
    # whose bucket center lands on our 15-character-wide target.
    $ws.Columns.Item(5).ColumnWidth = 14.166666666666666
}
